# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets,
# matching the gh-pages output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new value for column F
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 644
    3  = 6086
    12 = 1990
    19 = 372
    21 = 103
    23 = 1191
    24 = 2983
    26 = 2534
    27 = 4329
    31 = 1368
    34 = 49
    37 = 1048
    38 = 1301
    40 = 1133
    49 = 3627
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 644
    3  = 6086
    10 = 1990
    19 = 372
    20 = 103
    23 = 1191
    25 = 2983
    26 = 2534
    27 = 4329
    31 = 1368
    33 = 1048
    35 = 1301
    37 = 1133
    48 = 3627
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}

$wb.Save()
